$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$block1 = New-Object 'object[,]' 24,5
$block1[0,0] = 1.02
$block1[0,1] = 1.04771439876932
$block1[0,2] = 1.049524829611311
$block1[0,3] = 1.061309456830279
$block1[0,4] = 1.068673395832178
$block1[1,0] = 1.02
$block1[1,1] = 1.048671758759837
$block1[1,2] = 1.050257835506819
$block1[1,3] = 1.06225734506084
$block1[1,4] = 1.069700470109035
$block1[2,0] = 1.02
$block1[2,1] = 1.049291704404007
$block1[2,2] = 1.050732484695203
$block1[2,3] = 1.062871530110706
$block1[2,4] = 1.07036604001261
$block1[3,0] = 1.02
$block1[3,1] = 1.049552441153311
$block1[3,2] = 1.050932108670547
$block1[3,3] = 1.063129932972769
$block1[3,4] = 1.070646080094037
$block1[4,0] = 1.02
$block1[4,1] = 1.049596226520543
$block1[4,2] = 1.050965631155504
$block1[4,3] = 1.063173331613907
$block1[4,4] = 1.070693113774386
$block1[5,0] = 1.02
$block1[5,1] = 1.049295187943749
$block1[5,2] = 1.050735151761995
$block1[5,3] = 1.062874982120057
$block1[5,4] = 1.070369781001874
$block1[6,0] = 1.02
$block1[6,1] = 1.048037845290714
$block1[6,2] = 1.049772480141166
$block1[6,3] = 1.061629625895582
$block1[6,4] = 1.069020295941049
$block1[7,0] = 1.02
$block1[7,1] = 1.045825881879337
$block1[7,2] = 1.048078828779543
$block1[7,3] = 1.059441617726245
$block1[7,4] = 1.066649917788897
$block1[8,0] = 1.02
$block1[8,1] = 1.044353732791398
$block1[8,2] = 1.04695161466408
$block1[8,3] = 1.057987362029752
$block1[8,4] = 1.06507483213383
$block1[9,0] = 1.02
$block1[9,1] = 1.043716877251573
$block1[9,2] = 1.046463981018691
$block1[9,3] = 1.057358714624933
$block1[9,4] = 1.064394041478182
$block1[10,0] = 1.02
$block1[10,1] = 1.043480410839417
$block1[10,2] = 1.046282922087442
$block1[10,3] = 1.057125366554235
$block1[10,4] = 1.064141351621544
$block1[11,0] = 1.02
$block1[11,1] = 1.043531129598765
$block1[11,2] = 1.046321756659967
$block1[11,3] = 1.05717541326948
$block1[11,4] = 1.064195546015226
$block1[12,0] = 1.02
$block1[12,1] = 1.043697329022039
$block1[12,2] = 1.04644901319784
$block1[12,3] = 1.05733942273927
$block1[12,4] = 1.064373150253303
$block1[13,0] = 1.02
$block1[13,1] = 1.043799741932209
$block1[13,2] = 1.046527429450794
$block1[13,3] = 1.057440495562949
$block1[13,4] = 1.064482602777891
$block1[14,0] = 1.02
$block1[14,1] = 1.044396011389436
$block1[14,2] = 1.046983987043093
$block1[14,3] = 1.05802910560525
$block1[14,4] = 1.065120040009108
$block1[15,0] = 1.02
$block1[15,1] = 1.044770195213323
$block1[15,2] = 1.047270496879738
$block1[15,3] = 1.058398608322805
$block1[15,4] = 1.065520218069164
$block1[16,0] = 1.02
$block1[16,1] = 1.044988507639993
$block1[16,2] = 1.047437657283436
$block1[16,3] = 1.058614234717367
$block1[16,4] = 1.065753753944782
$block1[17,0] = 1.02
$block1[17,1] = 1.045062956238017
$block1[17,2] = 1.047494662099252
$block1[17,3] = 1.058687774964814
$block1[17,4] = 1.065833403759545
$block1[18,0] = 1.02
$block1[18,1] = 1.044730042894665
$block1[18,2] = 1.047239752519536
$block1[18,3] = 1.058358953649898
$block1[18,4] = 1.065477270456691
$block1[19,0] = 1.02
$block1[19,1] = 1.043648384941364
$block1[19,2] = 1.046411537379692
$block1[19,3] = 1.05729112162477
$block1[19,4] = 1.064320845086097
$block1[20,0] = 1.02
$block1[20,1] = 1.042968825266626
$block1[20,2] = 1.045891210750366
$block1[20,3] = 1.056620656984071
$block1[20,4] = 1.063594831820454
$block1[21,0] = 1.02
$block1[21,1] = 1.043329023016522
$block1[21,2] = 1.046167006947909
$block1[21,3] = 1.05697599505544
$block1[21,4] = 1.063979602729641
$block1[22,0] = 1.02
$block1[22,1] = 1.044748185822005
$block1[22,2] = 1.047253644436503
$block1[22,3] = 1.058376871575484
$block1[22,4] = 1.065496676267159
$block1[23,0] = 1.02
$block1[23,1] = 1.04639729162483
$block1[23,2] = 1.048516350974731
$block1[23,3] = 1.060006496798485
$block1[23,4] = 1.067261811542108
$ws.Range("B2:F25").Value = $block1

$block2 = New-Object 'object[,]' 24,6
$block2[0,0] = 1.044722035198725
$block2[0,1] = 1.052761996675746
$block2[0,2] = 1.052281646575665
$block2[0,3] = 1.064033867021544
$block2[0,4] = 1.071377948351662
$block2[0,5] = 1.054257039464515
$block2[1,0] = 1.044977003681923
$block2[1,1] = 1.053367902734957
$block2[1,2] = 1.052827120638053
$block2[1,3] = 1.064796015085454
$block2[1,4] = 1.072220521714395
$block2[1,5] = 1.054863805979828
$block2[2,0] = 1.045141054786249
$block2[2,1] = 1.053759799607003
$block2[2,2] = 1.053179754618427
$block2[2,3] = 1.065289386379099
$block2[2,4] = 1.072766094227285
$block2[2,5] = 1.055256259390383
$block2[3,0] = 1.045209798613667
$block2[3,1] = 1.053924512753279
$block2[3,2] = 1.053327923468238
$block2[3,3] = 1.065496849107424
$block2[3,4] = 1.0729955407772
$block2[3,5] = 1.05542120644822
$block2[4,0] = 1.045221327898979
$block2[4,1] = 1.053952166428106
$block2[4,2] = 1.053352797043345
$block2[4,3] = 1.06553168588337
$block2[4,4] = 1.073034070985975
$block2[4,5] = 1.055448899394437
$block2[5,0] = 1.045141974222476
$block2[5,1] = 1.053762000671506
$block2[5,2] = 1.053181734767276
$block2[5,3] = 1.06529215831318
$block2[5,4] = 1.072769159759311
$block2[5,5] = 1.05525846358065
$block2[6,0] = 1.044808395379069
$block2[6,1] = 1.052966799101831
$block2[6,2] = 1.052466058702769
$block2[6,3] = 1.06429139435454
$block2[6,4] = 1.071662622500444
$block2[6,5] = 1.054462132733535
$block2[7,0] = 1.044213481169763
$block2[7,1] = 1.051564326324872
$block2[7,2] = 1.051202500860411
$block2[7,3] = 1.062529575473211
$block2[7,4] = 1.069715649734291
$block2[7,5] = 1.053057668284367
$block2[8,0] = 1.043812122498263
$block2[8,1] = 1.05062856878461
$block2[8,2] = 1.050358535387697
$block2[8,3] = 1.061356198895696
$block2[8,4] = 1.068419667776962
$block2[8,5] = 1.052120581861065
$block2[9,0] = 1.043637209120667
$block2[9,1] = 1.050223201437517
$block2[9,2] = 1.049992721446314
$block2[9,3] = 1.060848402853544
$block2[9,4] = 1.067858979656426
$block2[9,5] = 1.051714638845843
$block2[10,0] = 1.043572070307992
$block2[10,1] = 1.05007260395496
$block2[10,2] = 1.04985678687607
$block2[10,3] = 1.060659828320535
$block2[10,4] = 1.067650788256788
$block2[10,5] = 1.051563827497587
$block2[11,0] = 1.043586050408642
$block2[11,1] = 1.050104908790461
$block2[11,2] = 1.049885947778087
$block2[11,3] = 1.060700276213917
$block2[11,4] = 1.067695442694412
$block2[11,5] = 1.05159617820966
$block2[12,0] = 1.043631828154714
$block2[12,1] = 1.050210753528603
$block2[12,2] = 1.049981486169351
$block2[12,3] = 1.060832814330698
$block2[12,4] = 1.067841768993128
$block2[12,5] = 1.05170217325947
$block2[13,0] = 1.04366001106037
$block2[13,1] = 1.050275964542663
$block2[13,2] = 1.050040343263431
$block2[13,3] = 1.060914481233299
$block2[13,4] = 1.067931935173726
$block2[13,5] = 1.051767476880651
$block2[14,0] = 1.043823707311972
$block2[14,1] = 1.050655467992079
$block2[14,2] = 1.050382805499316
$block2[14,3] = 1.061389905701429
$block2[14,4] = 1.068456889005929
$block2[14,5] = 1.052147519268495
$block2[15,0] = 1.043926089447364
$block2[15,1] = 1.0508934733212
$block2[15,2] = 1.05059752416442
$block2[15,3] = 1.061688203510962
$block2[15,4] = 1.068786308003226
$block2[15,5] = 1.052385862592482
$block2[16,0] = 1.043985698887746
$block2[16,1] = 1.0510322805619
$block2[16,2] = 1.050722730047866
$block2[16,3] = 1.061862222894889
$block2[16,4] = 1.068978498841528
$block2[16,5] = 1.052524866955385
$block2[17,0] = 1.044006005779128
$block2[17,1] = 1.051079607296259
$block2[17,2] = 1.050765415908065
$block2[17,3] = 1.061921563614242
$block2[17,4] = 1.069044038777705
$block2[17,5] = 1.052572260899136
$block2[18,0] = 1.043915116014789
$block2[18,1] = 1.050867939379616
$block2[18,2] = 1.050574490589273
$block2[18,3] = 1.061656196146895
$block2[18,4] = 1.068750959665579
$block2[18,5] = 1.052360292389773
$block2[19,0] = 1.043618352389004
$block2[19,1] = 1.050179585602923
$block2[19,2] = 1.049953354005743
$block2[19,3] = 1.060793783949977
$block2[19,4] = 1.067798677518963
$block2[19,5] = 1.051670961071761
$block2[20,0] = 1.043430792405216
$block2[20,1] = 1.049746640210758
$block2[20,2] = 1.049562503058102
$block2[20,3] = 1.060251803450188
$block2[20,4] = 1.067200363523307
$block2[20,5] = 1.05123740084748
$block2[21,0] = 1.043530313576204
$block2[21,1] = 1.049976166721682
$block2[21,2] = 1.049769730388688
$block2[21,3] = 1.060539093354985
$block2[21,4] = 1.067517500631775
$block2[21,5] = 1.051467253312376
$block2[22,0] = 1.043920074771366
$block2[22,1] = 1.050879477122047
$block2[22,2] = 1.050584898581344
$block2[22,3] = 1.06167065881263
$block2[22,4] = 1.06876693191478
$block2[22,5] = 1.052371846517121
$block2[23,0] = 1.044368119910328
$block2[23,1] = 1.051927039443792
$block2[23,2] = 1.051529445428036
$block2[23,3] = 1.062984845932385
$block2[23,4] = 1.070218640677447
$block2[23,5] = 1.05342089649752
$ws.Range("I2:N25").Value = $block2
